$wb = $excel.ActiveWorkbook

$wsOrg = $wb.Worksheets.Item("OrgData")
$wsOrg.Range("A2").Value = "AUTO_ORG_HCDAJ"

$wsSubOrg = $wb.Worksheets.Item("SubOrgData")
$wsSubOrg.Range("A2").Value = "AUTO_SUB_ORG_YMKWH"

$wsMember = $wb.Worksheets.Item("MemberData")
$wsMember.Range("A2").Value = "First9290"
